$d = $word.ActiveDocument

# The document has two "<id>...</id>" tag-id markers, each split across three
# separately-formatted runs: "<id>" (Courier New / gold), the bare id text
# (plain black), and "</id>" (Courier New / gold again). Newly downloaded
# tc/tcn/tl content consolidates each of these into a single run (using the
# "<id>" run's Courier-New/gold formatting) and, for the second marker,
# renames the id from "p160r_a2" to "p160r_2".

# 1) <id>p160r_1</id> -- text unchanged, just merge the 3 runs into 1.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("<id>p160r_1</id>", $true, $false, $false, $false, $false,
                              $true, 1, $false, "<id>p160r_1</id>", 2)
if (-not $found1) {
    throw "Could not find '<id>p160r_1</id>' to consolidate"
}

# 2) <id>p160r_a2</id> -- merge the 3 runs into 1 AND rename p160r_a2 -> p160r_2.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("<id>p160r_a2</id>", $true, $false, $false, $false, $false,
                              $true, 1, $false, "<id>p160r_2</id>", 2)
if (-not $found2) {
    throw "Could not find '<id>p160r_a2</id>' to consolidate/rename"
}

Write-Output "Consolidated id run 1: $found1; id run 2 (renamed to p160r_2): $found2"
